$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.276.44"
$ws.Range("E2").Value = "  +0.76%  "
$ws.Range("D3").Value = "3.907.00"
$ws.Range("E3").Value = "  +1.40%  "
$ws.Range("E4").Value = "  +0.05%  "
$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "469.08"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  +10.07%  "
$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.18"
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = "  +11.62%  "
$style = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.634"
$ws.Range("D7").Style = $style
$ws.Range("E7").Value = "  +3.73%  "
$style = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = $style
$ws.Range("E8").Value = "  +0.09%  "
$style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.743"
$ws.Range("D9").Style = $style
$ws.Range("E9").Value = "  +2.31%  "
$style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.157"
$ws.Range("D10").Style = $style
$style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000315"
$ws.Range("D11").Style = $style
$style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.53"
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = "  +5.89%  "
$ws.Range("D13").Value = "4.569.20"
$ws.Range("E13").Value = "  +2.06%  "
$style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.39"
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = "  +1.00%  "
$style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.86"
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = "  -7.24%  "
$ws.Range("D16").Value = "3.926.38"
$ws.Range("E16").Value = "  +1.01%  "
$ws.Range("E17").Value = "  -0.31%  "
$style = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.06"
$ws.Range("D18").Style = $style
$ws.Range("E18").Value = "  +0.39%  "
$style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.16"
$ws.Range("D19").Style = $style
$ws.Range("E19").Value = "  +7.16%  "
$ws.Range("D20").Value = "67.670.85"
$ws.Range("E20").Value = "  +0.96%  "
$style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "431.57"
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = "  +4.46%  "
$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.74"
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = "  -1.56%  "
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "88.89"
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = "  +5.34%  "
$ws.Range("B24").Value = "ImmutableX"
$ws.Range("C24").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.27"
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = "  +7.34%  "
$style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.60"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = "  +10.59%  "
$style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "38.04"
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = "  +0.84%  "
$style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.20"
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = "  +12.66%  "
$style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.04"
$ws.Range("D28").Style = $style
$style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.48"
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = "  +2.41%  "
$style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "731.10"
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = "  +1.67%  "
$ws.Range("E31").Value = "  +9.66%  "
$style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.72"
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = "  +1.15%  "
$style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.76"
$ws.Range("D33").Style = $style
$ws.Range("E33").Value = "  -0.68%  "
$style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "43.36"
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = "  +10.51%  "
$ws.Range("E35").Value = "  +7.47%  "
$style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.41"
$ws.Range("D36").Style = $style
$ws.Range("E36").Value = "  +2.78%  "
$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.49"
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = "  +2.68%  "
$style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0481"
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = "  +4.03%  "
$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.346"
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = "  +10.19%  "
$ws.Range("B41").Value = "ThetaToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.91"
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = "  +0.69%  "
$ws.Range("B42").Value = "PEPE"
$ws.Range("C42").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D42").Value = "0.0₃0688"
$ws.Range("E42").Value = "  -8.88%  "
$style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.141"
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = "  +4.00%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("D44").Style = $style
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("B45").Value = "LidoDAOToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.49"
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = "  +3.45%  "
$ws.Range("E46").Value = "  +11.53%  "
$style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.78"
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = "  +7.23%  "
$style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.27"
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = "  +1.81%  "
$style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.18"
$ws.Range("D49").Style = $style
$ws.Range("E49").Value = "  +6.45%  "
$style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "144.26"
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = "  +1.39%  "
$style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.88"
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = "  +2.09%  "
